$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Start from a clean slate for cell contents (keeps column defs which we'll refit below)
$ws.Cells.Clear()
$ws.Columns("C").ClearFormats()
$ws.Columns("D").ClearFormats()

# Header row
$ws.Range("A1").Value = "name"
$ws.Range("B1").Value = "platform"
$ws.Range("C1").Value = "code"
$ws.Range("D1").Value = "cookie"
$ws.Range("E1").Value = "email"
$ws.Range("F1").Value = "phone"
$ws.Range("G1").Value = "rate"

$ws.Range("H1").NumberFormat = "@"
$ws.Range("H1").Value = "merchant_id"

# Row 2 - Shopee store
$ws.Range("A2").Value = "storexeonvn"
$ws.Range("B2").Value = "shopee"
$ws.Range("C2").Formula = "=B2&`":`"&A2"
$ws.Range("G2").Value = 1.05

$ws.Range("H2").NumberFormat = "@"
$ws.Range("H2").Value = "76922911"

# Row 3 - Lazada store 1
$ws.Range("A3").Value = "thexeonstore"
$ws.Range("B3").Value = "lazada"
$ws.Range("C3").Formula = "=B3&`":`"&A3"
$ws.Range("G3").Value = 1.05

$ws.Range("H3").NumberFormat = "@"
$ws.Range("H3").Value = "0"

# Row 4 - Lazada store 2
$ws.Range("A4").Value = "pczone"
$ws.Range("B4").Value = "lazada"
$ws.Range("C4").Formula = "=B4&`":`"&A4"
$ws.Range("G4").Value = 1.05

$ws.Range("H4").NumberFormat = "@"
$ws.Range("H4").Value = "0"

# Column widths to match bestFit behavior
$ws.Columns("A").AutoFit()
$ws.Columns("C").AutoFit()
$ws.Columns("G").AutoFit()
$ws.Columns("H").AutoFit()

$null = $ws.Range("D11").Select()
